# "contingencies with rene fine"
# Insert line7 / line8 into the lines table (pushing the extr* rows down by
# two) and append two more extr rows (extr7, extr8) at the bottom, giving
# every row fresh from_bus/to_bus/in_service data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: was extr1, becomes line7 -------------------------------------
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- Row 9: was extr2, becomes line8 -------------------------------------
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Row 10: was extr3, becomes extr1 -------------------------------------
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Row 11: was extr4, becomes extr2 -------------------------------------
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# --- Row 12: was extr5, becomes extr3 -------------------------------------
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# --- Row 13: was extr6, becomes extr4 -------------------------------------
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# --- Row 14: was extr7, becomes extr5 -------------------------------------
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# --- Row 15: was extr8, becomes extr6 -------------------------------------
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Row 16: new, extr7 ----------------------------------------------------
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# --- Row 17: new, extr8 ----------------------------------------------------
$ws.Range("A15").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
